$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.08%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'38.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'8.19%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.123"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.50%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.08172"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.63%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'2.024"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'8.77%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'7.916"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.12%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.9324"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.34%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1400"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'4.12%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1956"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.74%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.09198"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.25%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.03438"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.07%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.09860"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.30%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001412"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.78%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.006252"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'4.07%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'3.621"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-3.14%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'4.194"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.95%"
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'1.65%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.3450"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.21%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1314"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.93%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.813"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-6.63%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.2452"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.57%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04478"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.65%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001242"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.36%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D27").Value = "'0.0001303"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.37%"
$ws.Range("E27").Style = "Normal"

$ws.Range("D39").Value = "'0.02132"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'9.99%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.05190"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-2.22%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007468"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.80%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.01001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.87%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.1366"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.31%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.002135"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-0.56%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.009763"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-3.93%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00006338"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.90%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.39%"
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'-0.63%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.001603"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-3.21%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.00002105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.39%"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0002005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.39%"
$ws.Range("E51").Style = "Normal"
